# Update cryptocurrency Price (D) and Volume(1h) (E) columns to the
# latest scraped figures, per the 2023-07-14 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Range("D2")
$cD.NumberFormat = "@"
$cD.Value = '31.469.61'
$cD.Style = "Normal"
$ws.Range("E2").Value = '  +3.82%  '

$cD = $ws.Range("D3")
$cD.NumberFormat = "@"
$cD.Value = '2.008.70'
$cD.Style = "Normal"
$ws.Range("E3").Value = '  +7.54%  '

$ws.Range("E4").Value = '  -0.07%  '

$cD = $ws.Range("D5")
$cD.NumberFormat = "@"
$cD.Value = '0.7699'
$cD.Style = "Normal"
$ws.Range("E5").Value = '  +63.16%  '

$cD = $ws.Range("D6")
$cD.NumberFormat = "@"
$cD.Value = '259.48'
$cD.Style = "Normal"
$ws.Range("E6").Value = '  +6.68%  '

$cD = $ws.Range("D7")
$cD.NumberFormat = "@"
$cD.Value = '0.9996'
$cD.Style = "Normal"
$ws.Range("E7").Value = '  -0.11%  '

$cD = $ws.Range("D8")
$cD.NumberFormat = "@"
$cD.Value = '0.3583'
$cD.Style = "Normal"
$ws.Range("E8").Value = '  +24.92%  '

$cD = $ws.Range("D9")
$cD.NumberFormat = "@"
$cD.Value = '28.34'
$cD.Style = "Normal"
$ws.Range("E9").Value = '  +31.11%  '

$cD = $ws.Range("D10")
$cD.NumberFormat = "@"
$cD.Value = '0.07072'
$cD.Style = "Normal"
$ws.Range("E10").Value = '  +9.29%  '

$cD = $ws.Range("D11")
$cD.NumberFormat = "@"
$cD.Value = '0.8422'
$cD.Style = "Normal"
$ws.Range("E11").Value = '  +17.85%  '

$cD = $ws.Range("D12")
$cD.NumberFormat = "@"
$cD.Value = '0.08109'
$cD.Style = "Normal"
$ws.Range("E12").Value = '  +4.28%  '

$cD = $ws.Range("D13")
$cD.NumberFormat = "@"
$cD.Value = '2.006.38'
$cD.Style = "Normal"
$ws.Range("E13").Value = '  +7.42%  '

$cD = $ws.Range("D14")
$cD.NumberFormat = "@"
$cD.Value = '101.41'
$cD.Style = "Normal"
$ws.Range("E14").Value = '  +5.76%  '

$cD = $ws.Range("D15")
$cD.NumberFormat = "@"
$cD.Value = '5.635'
$cD.Style = "Normal"
$ws.Range("E15").Value = '  +10.20%  '

$cD = $ws.Range("D16")
$cD.NumberFormat = "@"
$cD.Value = '275.07'
$cD.Style = "Normal"
$ws.Range("E16").Value = '  -2.78%  '

$cD = $ws.Range("D17")
$cD.NumberFormat = "@"
$cD.Value = '31.477.45'
$cD.Style = "Normal"
$ws.Range("E17").Value = '  +3.89%  '

$cD = $ws.Range("D18")
$cD.NumberFormat = "@"
$cD.Value = '14.69'
$cD.Style = "Normal"
$ws.Range("E18").Value = '  +13.26%  '

$cD = $ws.Range("D19")
$cD.NumberFormat = "@"
$cD.Value = '5.941'
$cD.Style = "Normal"
$ws.Range("E19").Value = '  +13.54%  '

$cD = $ws.Range("D20")
$cD.NumberFormat = "@"
$cD.Value = '0.000007976'
$cD.Style = "Normal"
$ws.Range("E20").Value = '  +7.05%  '

$cD = $ws.Range("D21")
$cD.NumberFormat = "@"
$cD.Value = '2.272.26'
$cD.Style = "Normal"
$ws.Range("E21").Value = '  +7.82%  '

$cD = $ws.Range("D22")
$cD.NumberFormat = "@"
$cD.Value = '0.9993'
$cD.Style = "Normal"
$ws.Range("E22").Value = '  -0.12%  '

$cD = $ws.Range("D23")
$cD.NumberFormat = "@"
$cD.Value = '0.9995'
$cD.Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '

$cD = $ws.Range("D24")
$cD.NumberFormat = "@"
$cD.Value = '7.228'
$cD.Style = "Normal"
$ws.Range("E24").Value = '  +15.88%  '

$ws.Range("E25").Value = '  +13.25%  '

$cD = $ws.Range("D26")
$cD.NumberFormat = "@"
$cD.Value = '164.08'
$cD.Style = "Normal"
$ws.Range("E26").Value = '  +0.96%  '

$cD = $ws.Range("D27")
$cD.NumberFormat = "@"
$cD.Value = '0.1460'
$cD.Style = "Normal"
$ws.Range("E27").Value = '  +52.19%  '

$cD = $ws.Range("D28")
$cD.NumberFormat = "@"
$cD.Value = '20.15'
$cD.Style = "Normal"
$ws.Range("E28").Value = '  +7.89%  '

$ws.Range("E29").Value = '  +27.50%  '

$cD = $ws.Range("D30")
$cD.NumberFormat = "@"
$cD.Value = '1.629'
$cD.Style = "Normal"
$ws.Range("E30").Value = '  +10.00%  '

$cD = $ws.Range("D31")
$cD.NumberFormat = "@"
$cD.Value = '4.637'
$cD.Style = "Normal"
$ws.Range("E31").Value = '  +10.51%  '

$ws.Range("E32").Value = '  +3.41%  '

$cD = $ws.Range("D33")
$cD.NumberFormat = "@"
$cD.Value = '4.407'
$cD.Style = "Normal"
$ws.Range("E33").Value = '  +7.25%  '

$cD = $ws.Range("D34")
$cD.NumberFormat = "@"
$cD.Value = '0.05215'
$cD.Style = "Normal"
$ws.Range("E34").Value = '  +8.49%  '

$cD = $ws.Range("D35")
$cD.NumberFormat = "@"
$cD.Value = '1.237'
$cD.Style = "Normal"
$ws.Range("E35").Value = '  +10.87%  '

$cD = $ws.Range("D36")
$cD.NumberFormat = "@"
$cD.Value = '0.7628'
$cD.Style = "Normal"
$ws.Range("E36").Value = '  +11.67%  '

$cD = $ws.Range("D37")
$cD.NumberFormat = "@"
$cD.Value = '2.800'
$cD.Style = "Normal"

$cD = $ws.Range("D38")
$cD.NumberFormat = "@"
$cD.Value = '0.02024'
$cD.Style = "Normal"

$cD = $ws.Range("D39")
$cD.NumberFormat = "@"
$cD.Value = '2.953'
$cD.Style = "Normal"
$ws.Range("E39").Value = '  +3.80%  '

$cD = $ws.Range("D40")
$cD.NumberFormat = "@"
$cD.Value = '6.723'
$cD.Style = "Normal"
$ws.Range("E40").Value = '  +8.30%  '

$cD = $ws.Range("D41")
$cD.NumberFormat = "@"
$cD.Value = '80.46'
$cD.Style = "Normal"
$ws.Range("E41").Value = '  +6.87%  '

$ws.Range("E42").Value = '  +14.30%  '

$cD = $ws.Range("D43")
$cD.NumberFormat = "@"
$cD.Value = '0.4755'
$cD.Style = "Normal"
$ws.Range("E43").Value = '  +13.71%  '

$cD = $ws.Range("D44")
$cD.NumberFormat = "@"
$cD.Value = '0.8641'
$cD.Style = "Normal"
$ws.Range("E44").Value = '  +5.08%  '

$cD = $ws.Range("D45")
$cD.NumberFormat = "@"
$cD.Value = '104.80'
$cD.Style = "Normal"
$ws.Range("E45").Value = '  +4.30%  '

$ws.Range("E46").Value = '  +0.03%  '

$cD = $ws.Range("D47")
$cD.NumberFormat = "@"
$cD.Value = '9.965'
$cD.Style = "Normal"
$ws.Range("E47").Value = '  +4.07%  '

$cD = $ws.Range("D48")
$cD.NumberFormat = "@"
$cD.Value = '7.665'
$cD.Style = "Normal"
$ws.Range("E48").Value = '  +9.76%  '

$cD = $ws.Range("D49")
$cD.NumberFormat = "@"
$cD.Value = '0.4374'
$cD.Style = "Normal"
$ws.Range("E49").Value = '  +13.01%  '

$cD = $ws.Range("D50")
$cD.NumberFormat = "@"
$cD.Value = '37.14'
$cD.Style = "Normal"
$ws.Range("E50").Value = '  +6.23%  '

$cD = $ws.Range("D51")
$cD.NumberFormat = "@"
$cD.Value = '944.16'
$cD.Style = "Normal"
$ws.Range("E51").Value = '  +5.74%  '
